$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 3.272327238179451
$ws.Range("C2").Value = 1.626987699542094
$ws.Range("D2").Value = 3.223369029078222
$ws.Range("E2").Value = 2797.565817734744
$ws.Range("G2").Value = 2805.688501701543
